$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "Globo"
$ws.Range("B39").Value = "Bom Dia Inter"
$ws.Range("C39").Value = "Limpeza Pública"
$ws.Range("D39").Value = "2025-04-02T11:14"
$ws.Range("E39").Value = "Positivo"
$ws.Range("F39").Value = "Teste 2"
